# "Add files via upload" — re-upload of meta/en/3-4-1.xlsx with an updated
# organization website and the selection left parked on that cell.
#
# Semantic edit performed by the author:
#   - Row 10 ("Organization website (if available)") value changed from
#     "www.stat.kg " to "www.stat.gov.kg".
#   - The sheet's active selection moved from B2 to B10 (the cell that was
#     just edited) before saving.
#
# Setting the cell's .Value to a brand-new string naturally drops the old,
# now-unreferenced shared string and appends the new one at the end of the
# shared-strings table — which is exactly the reordering/renumbering seen
# across B10..B26 in the target diff (every t="s" index from the old
# "3. Definitions..." block onward shifts down by one once the stale
# "www.stat.kg " entry is garbage-collected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website cell.
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move/leave the selection on the edited cell.
$ws.Range("B10").Select()
